# Applies the "Updated cryptos list" data refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that assigns a literal text value to a cell without Excel
# reinterpreting numeric-looking strings (e.g. "246.48") as numbers,
# while restoring the cells original style so no formatting drifts.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '37.333.08'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '2.023.31'
$ws.Range('E3').Value = '  +2.82%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '246.48'
$ws.Range('E5').Value = '  +1.09%  '
Set-TextValue $ws.Range('D6') '0.619'
$ws.Range('E6').Value = '  -1.08%  '
Set-TextValue $ws.Range('D7') '57.91'
$ws.Range('E7').Value = '  -3.42%  '
$ws.Range('E8').Value = '  +0.09%  '
Set-TextValue $ws.Range('D9') '0.387'
$ws.Range('E9').Value = '  +2.31%  '
Set-TextValue $ws.Range('D10') '0.0799'
$ws.Range('E10').Value = '  +1.29%  '
Set-TextValue $ws.Range('D11') '0.103'
$ws.Range('E11').Value = '  -0.30%  '
Set-TextValue $ws.Range('D12') '14.91'
$ws.Range('E12').Value = '  +4.81%  '
$ws.Range('D13').Value = '2.328.75'
$ws.Range('E13').Value = '  +3.17%  '
Set-TextValue $ws.Range('D14') '0.833'
$ws.Range('E14').Value = '  -0.84%  '
Set-TextValue $ws.Range('D15') '21.51'
$ws.Range('E15').Value = '  -0.38%  '
Set-TextValue $ws.Range('D16') '5.37'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').Value = '2.025.24'
$ws.Range('E17').Value = '  +2.88%  '
$ws.Range('D18').Value = '37.279.96'
$ws.Range('E18').Value = '  +2.06%  '
Set-TextValue $ws.Range('D19') '69.92'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '0.0₃0852'
$ws.Range('E20').Value = '  -0.18%  '
Set-TextValue $ws.Range('D21') '5.19'
$ws.Range('E21').Value = '  +2.33%  '
Set-TextValue $ws.Range('D22') '227.26'
$ws.Range('E22').Value = '  -0.18%  '
Set-TextValue $ws.Range('D23') '0.999'
$ws.Range('E23').Value = '  +0.06%  '
Set-TextValue $ws.Range('D24') '2.53'
$ws.Range('E24').Value = '  +4.36%  '
$ws.Range('E25').Value = '  -0.59%  '
Set-TextValue $ws.Range('D26') '9.15'
$ws.Range('E26').Value = '  +0.42%  '
Set-TextValue $ws.Range('D27') '163.32'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -5.48%  '
Set-TextValue $ws.Range('D29') '19.74'
$ws.Range('E29').Value = '  +1.80%  '
Set-TextValue $ws.Range('D30') '1.35'
$ws.Range('E30').Value = '  +0.21%  '
Set-TextValue $ws.Range('D31') '0.120'
$ws.Range('E31').Value = '  -0.49%  '
Set-TextValue $ws.Range('D32') '4.75'
$ws.Range('E32').Value = '  -0.88%  '
Set-TextValue $ws.Range('D33') '0.0665'
$ws.Range('E33').Value = '  +8.43%  '
Set-TextValue $ws.Range('D34') '4.56'
$ws.Range('E34').Value = '  +0.90%  '
Set-TextValue $ws.Range('D35') '2.47'
$ws.Range('E35').Value = '  +8.01%  '
Set-TextValue $ws.Range('D36') '3.54'
$ws.Range('E36').Value = '  +5.75%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +2.27%  '
Set-TextValue $ws.Range('D39') '5.38'
$ws.Range('E39').Value = '  -0.41%  '
Set-TextValue $ws.Range('D40') '3.00'
$ws.Range('E40').Value = '  +3.18%  '
Set-TextValue $ws.Range('D41') '0.0968'
$ws.Range('E41').Value = '  +0.30%  '
Set-TextValue $ws.Range('D42') '0.0217'
$ws.Range('E42').Value = '  +3.57%  '
Set-TextValue $ws.Range('D43') '1.16'
$ws.Range('E43').Value = '  -0.32%  '
Set-TextValue $ws.Range('D44') '16.37'
$ws.Range('E44').Value = '  +3.14%  '
$ws.Range('D45').Value = '1.391.22'
$ws.Range('E45').Value = '  +1.87%  '
Set-TextValue $ws.Range('D46') '90.56'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D47') '7.45'
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D48') '1.04'
$ws.Range('E48').Value = '  +1.81%  '
Set-TextValue $ws.Range('D49') '2.08'
$ws.Range('E49').Value = '  +12.59%  '
Set-TextValue $ws.Range('D50') '2.87'
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('D51').Value = '2.219.97'
$ws.Range('E51').Value = '  +3.20%  '
